$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - update column F (想去人数) values
$wsExh = $wb.Worksheets.Item("展览")
$wsExh.Range("F4").Value = 640
$wsExh.Range("F5").Value = 198
$wsExh.Range("F6").Value = 1
$wsExh.Range("F7").Value = 9671
$wsExh.Range("F10").Value = 1224
$wsExh.Range("F11").Value = 2781
$wsExh.Range("F12").Value = 164
$wsExh.Range("F13").Value = 108
$wsExh.Range("F14").Value = 16
$wsExh.Range("F15").Value = 25
$wsExh.Range("F17").Value = 489
$wsExh.Range("F18").Value = 101
$wsExh.Range("F20").Value = 1375

# Sheet "全部类型" (All Types) - update column F (想去人数) values
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 640
$wsAll.Range("F6").Value = 198
$wsAll.Range("F7").Value = 1
$wsAll.Range("F8").Value = 9671
$wsAll.Range("F11").Value = 1224
$wsAll.Range("F12").Value = 2781
$wsAll.Range("F13").Value = 164
$wsAll.Range("F14").Value = 108
$wsAll.Range("F15").Value = 16
$wsAll.Range("F16").Value = 25
$wsAll.Range("F18").Value = 489
$wsAll.Range("F19").Value = 101
$wsAll.Range("F21").Value = 1375
